$d = $word.ActiveDocument
$d.Content.Find.Execute("Images:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Images: IconFinder, Unsplash, CloudConvert png to ico", 2)
